$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "65.680.78"
$ws.Cells.Item(2, 5).Value = "  +0.67%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.589.11"
$ws.Cells.Item(3, 5).Value = "  +1.46%  "
$ws.Cells.Item(4, 5).Value = "  +0.05%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "603.40"
$ws.Cells.Item(5, 5).Value = "  +0.64%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "137.58"
$ws.Cells.Item(6, 5).Value = "  -1.06%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.588.77"
$ws.Cells.Item(7, 5).Value = "  +1.44%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  +1.48%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.126"
$ws.Cells.Item(10, 5).Value = "  +0.93%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.25"
$ws.Cells.Item(11, 5).Value = "  +4.63%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.392"
$ws.Cells.Item(12, 5).Value = "  +0.54%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.196.96"
$ws.Cells.Item(13, 5).Value = "  +1.63%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "28.17"
$ws.Cells.Item(14, 5).Value = "  +3.72%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.0000187"
$ws.Cells.Item(15, 5).Value = "  +1.25%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.587.35"
$ws.Cells.Item(16, 5).Value = "  +1.25%  "
$ws.Cells.Item(17, 5).Value = "  -0.11%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "65.792.24"
$ws.Cells.Item(18, 5).Value = "  +0.80%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.03"
$ws.Cells.Item(19, 5).Value = "  -2.96%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.65"
$ws.Cells.Item(20, 5).Value = "  +2.50%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.90"
$ws.Cells.Item(21, 5).Value = "  -1.09%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "397.03"
$ws.Cells.Item(22, 5).Value = "  +0.85%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.591"
$ws.Cells.Item(23, 5).Value = "  +3.11%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "3.733.88"
$ws.Cells.Item(24, 5).Value = "  +1.76%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "74.30"
$ws.Cells.Item(25, 5).Value = "  +0.68%  "
$ws.Cells.Item(26, 5).Value = "  -0.03%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.0000119"
$ws.Cells.Item(27, 5).Value = "  +3.17%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.13"
$ws.Cells.Item(28, 5).Value = "  +5.85%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.67"
$ws.Cells.Item(29, 5).Value = "  +31.54%  "
$ws.Cells.Item(30, 5).Value = "  +5.69%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.62"
$ws.Cells.Item(31, 5).Value = "  +4.75%  "
$ws.Cells.Item(32, 5).Value = "  -1.64%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.594.49"
$ws.Cells.Item(33, 5).Value = "  +1.44%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "24.54"
$ws.Cells.Item(34, 5).Value = "  +3.15%  "
$ws.Cells.Item(35, 5).Value = "  -0.01%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.148"
$ws.Cells.Item(36, 5).Value = "  +1.70%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.41"
$ws.Cells.Item(37, 5).Value = "  +8.83%  "
$ws.Cells.Item(38, 5).Value = "  +4.11%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "7.09"
$ws.Cells.Item(39, 5).Value = "  +1.97%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "168.49"
$ws.Cells.Item(40, 5).Value = "  -0.48%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0839"
$ws.Cells.Item(41, 5).Value = "  +4.81%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.843"
$ws.Cells.Item(42, 5).Value = "  +2.14%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "26.73"
$ws.Cells.Item(43, 5).Value = "  +1.42%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.28"
$ws.Cells.Item(44, 5).Value = "  +8.41%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "43.18"
$ws.Cells.Item(45, 5).Value = "  +0.83%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "4.55"
$ws.Cells.Item(46, 5).Value = "  +2.80%  "
$ws.Cells.Item(47, 5).Value = "  +0.07%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.70"
$ws.Cells.Item(48, 5).Value = "  +1.65%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.05"
$ws.Cells.Item(49, 5).Value = "  +3.54%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.449.92"
$ws.Cells.Item(50, 5).Value = "  +1.82%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "318.97"
$ws.Cells.Item(51, 5).Value = "  +5.45%  "
